# Update the wallets worksheet: replace the old single destination/private-key
# pair with two new private-key / amount rows, and drop the now-unused
# "Destination" column values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Amount cells (B2, B3) ---------------------------------------------
# These must end up as *text* "0.001" (shared string), not as the number
# 0.001, and B3 must not pick up any extra number-format/quote-prefix
# styling. Typing a numeric-looking value into Value/Formula always gets
# parsed back into a number by this engine, and forcing text via an
# apostrophe or NumberFormat="@" bakes a new style (quotePrefix / "@")
# onto the cell. To avoid corrupting styles, stage the text value on a
# scratch cell far outside the used range, copy it, and paste-special
# (values only) into the real targets - this carries over the *value*
# without touching the destination cell's existing style.
$ws.Range("Z1").Formula = "'0.001"
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("B3").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("Z1").Clear()

# --- Private key cells (A2, A3) -----------------------------------------
# Hex strings like these are not number-like, so they stay text naturally.
$ws.Range("A2").Value = "0x5af2c455c5889333b321c6f55ff076ae56458d0f23ca4a6ded84b48a8e08c33a"
$ws.Range("A3").Value = "0x77f97d0a55d2fa5425b4e0779f55ad41f55e202e7b85b2dde1fdff93f4380133"

# --- Flag columns (C3, D3) ------------------------------------------------
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0

# --- Drop the old Destination values (column E no longer used) ----------
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()

# --- Selection moves to A3 ------------------------------------------------
$ws.Range("A3").Select() | Out-Null
